$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 30

# Text columns: force text storage so values like dates/leading-zero
# week numbers are not auto-converted to dates/numbers by Excel, then
# clear the temporary text number-format again so the new row keeps the
# same (default/no) cell style as every other data row.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2024-01-07"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "11:46:30"
$ws.Cells.Item($row, 2).ClearFormats()

$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 3).Value = "Sunday"
$ws.Cells.Item($row, 3).ClearFormats()

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "01"
$ws.Cells.Item($row, 4).ClearFormats()

# Numeric columns
$ws.Cells.Item($row, 5).Value = 140580
$ws.Cells.Item($row, 6).Value = 143107
$ws.Cells.Item($row, 7).Value = 172010
$ws.Cells.Item($row, 8).Value = 147324
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 118125
$ws.Cells.Item($row, 11).Value = 224635
$ws.Cells.Item($row, 12).Value = 249314
$ws.Cells.Item($row, 13).Value = 185218
$ws.Cells.Item($row, 14).Value = 110432
$ws.Cells.Item($row, 15).Value = 40644
$ws.Cells.Item($row, 16).Value = 30820
$ws.Cells.Item($row, 17).Value = 72520
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42157
$ws.Cells.Item($row, 20).Value = -1
